$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values replacing the old Strike# derived values in column G
$kValues = @{
    2 = 5
    3 = 5
    4 = 7
    5 = 6
    6 = 6
    7 = 2
    8 = 10
    9 = 5
    10 = 1
    11 = 3
    12 = 1
    13 = 3
    14 = 7
    15 = 6
    16 = 5
    17 = 7
    18 = 4
    19 = 5
    20 = 4
    21 = 2
    22 = 7
    23 = 5
    24 = 2
    25 = 4
    26 = 2
    27 = 3
    28 = 5
    29 = 1
    30 = 7
    31 = 4
    32 = 4
    33 = 5
    34 = 4
    35 = 2
    36 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
